# Applies the "Updated cryptos list" GitHub Actions price/volume refresh
# described by the commit. Each coin's Price (column D) and Volume(1h)
# (column E) text is refreshed in place; rows 45/46 additionally swap the
# Quant / PaxDollar entries (name, link, price, volume all change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) store plain text values (e.g. "30.338.78",
# "  +0.77%  ") in the workbook. Force the cells to stay text so Excel does not
# reinterpret them as numbers/dates (which would corrupt values like "1.000")
# when we assign the new strings below.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range('D2').Value = '30.338.78'
$ws.Range('E2').Value = '  +0.77%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = '1.869.78'
$ws.Range('E3').Value = '  +0.39%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  -0.18%  '

# Row 5 - BNB
$ws.Range('D5').Value = '235.95'
$ws.Range('E5').Value = '  +0.62%  '

# Row 6 - USDC
$ws.Range('E6').Value = '  -0.25%  '

# Row 7 - XRP
$ws.Range('E7').Value = '  +0.09%  '

# Row 8 - Cardano
$ws.Range('D8').Value = '0.2852'
$ws.Range('E8').Value = '  +0.82%  '

# Row 9 - Dogecoin
$ws.Range('D9').Value = '0.06557'
$ws.Range('E9').Value = '  -0.28%  '

# Row 10 - Solana
$ws.Range('D10').Value = '21.49'
$ws.Range('E10').Value = '  +5.82%  '

# Row 11 - TRON
$ws.Range('D11').Value = '0.07891'
$ws.Range('E11').Value = '  +1.45%  '

# Row 12 - Litecoin
$ws.Range('D12').Value = '98.22'
$ws.Range('E12').Value = '  +1.72%  '

# Row 13 - WrappedEther
$ws.Range('D13').Value = '1.880.38'
$ws.Range('E13').Value = '  -0.01%  '

# Row 14 - Polkadot
$ws.Range('D14').Value = '5.112'
$ws.Range('E14').Value = '  +1.04%  '

# Row 15 - Polygon
$ws.Range('D15').Value = '0.6773'
$ws.Range('E15').Value = '  +1.25%  '

# Row 16 - BitcoinCash
$ws.Range('D16').Value = '278.86'
$ws.Range('E16').Value = '  -0.99%  '

# Row 17 - WrappedBTC
$ws.Range('D17').Value = '30.332.57'
$ws.Range('E17').Value = '  +0.63%  '

# Row 18 - Dai
$ws.Range('E18').Value = '  -0.17%  '

# Row 19 - Avalanche
$ws.Range('D19').Value = '12.74'
$ws.Range('E19').Value = '  +1.65%  '

# Row 20 - Uniswap
$ws.Range('D20').Value = '5.475'
$ws.Range('E20').Value = '  +2.05%  '

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range('D21').Value = '2.118.92'
$ws.Range('E21').Value = '  -0.55%  '

# Row 22 - ShibaInu
$ws.Range('D22').Value = '0.000007318'
$ws.Range('E22').Value = '  +1.07%  '

# Row 23 - BinanceUSD
$ws.Range('E23').Value = '  -0.13%  '

# Row 24 - Chainlink
$ws.Range('D24').Value = '6.160'
$ws.Range('E24').Value = '  +0.20%  '

# Row 25 - Monero
$ws.Range('D25').Value = '165.54'
$ws.Range('E25').Value = '  -0.81%  '

# Row 26 - Cosmos
$ws.Range('D26').Value = '9.161'
$ws.Range('E26').Value = '  -1.70%  '

# Row 27 - EthereumClassic
$ws.Range('D27').Value = '19.17'
$ws.Range('E27').Value = '  +0.78%  '

# Row 28 - LidoDAOToken
$ws.Range('E28').Value = '  -0.79%  '

# Row 29 - Toncoin
$ws.Range('D29').Value = '1.380'
$ws.Range('E29').Value = '  +0.53%  '

# Row 30 - Stellar
$ws.Range('D30').Value = '0.09661'
$ws.Range('E30').Value = '  +0.13%  '

# Row 31 - Filecoin
$ws.Range('D31').Value = '4.406'
$ws.Range('E31').Value = '  +0.73%  '

# Row 32 - PancakeSwap
$ws.Range('D32').Value = '1.476'
$ws.Range('E32').Value = '  +0.67%  '

# Row 33 - InternetComputer(DFINITY)
$ws.Range('D33').Value = '4.109'
$ws.Range('E33').Value = '  +0.37%  '

# Row 34 - Hedera
$ws.Range('D34').Value = '0.04723'
$ws.Range('E34').Value = '  +1.49%  '

# Row 35 - ARBITRUM
$ws.Range('D35').Value = '1.128'
$ws.Range('E35').Value = '  +3.86%  '

# Row 36 - ImmutableX
$ws.Range('D36').Value = '0.7086'
$ws.Range('E36').Value = '  +1.19%  '

# Row 37 - HuobiToken
$ws.Range('D37').Value = '2.726'
$ws.Range('E37').Value = '  +0.17%  '

# Row 38 - VeChain
$ws.Range('D38').Value = '0.01864'
$ws.Range('E38').Value = '  +0.41%  '

# Row 39 - FraxShare
$ws.Range('D39').Value = '6.339'
$ws.Range('E39').Value = '  -1.03%  '

# Row 40 - MXToken
$ws.Range('D40').Value = '2.535'
$ws.Range('E40').Value = '  +0.94%  '

# Row 41 - Aave
$ws.Range('D41').Value = '74.34'
$ws.Range('E41').Value = '  +4.11%  '

# Row 42 - RenderToken
$ws.Range('D42').Value = '1.957'
$ws.Range('E42').Value = '  +1.01%  '

# Row 43 - TrustWalletToken
$ws.Range('D43').Value = '0.8516'
$ws.Range('E43').Value = '  -0.90%  '

# Row 44 - TheSandbox
$ws.Range('D44').Value = '0.4194'
$ws.Range('E44').Value = '  +0.76%  '

# Row 45 - Quant
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  -0.24%  '

# Row 46 - PaxDollar
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '103.95'
$ws.Range('E46').Value = '  +0.96%  '

# Row 47 - Aptos
$ws.Range('D47').Value = '7.209'
$ws.Range('E47').Value = '  +0.35%  '

# Row 48 - EnergySwap
$ws.Range('D48').Value = '9.245'
$ws.Range('E48').Value = '  +1.85%  '

# Row 49 - Maker
$ws.Range('D49').Value = '939.35'
$ws.Range('E49').Value = '  -4.54%  '

# Row 50 - Elrond
$ws.Range('D50').Value = '34.27'
$ws.Range('E50').Value = '  +1.31%  '

# Row 51 - Algorand
$ws.Range('D51').Value = '0.1124'
$ws.Range('E51').Value = '  -1.61%  '
